$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 12 de Mayo de 2020 a las 02:05"
$ws.Range("B4").Value = 1385745
$ws.Range("C4").Value = 18107
$ws.Range("D4").Value = 260734
$ws.Range("E4").Value = 1043225
$ws.Range("F4").Value = 16484
$ws.Range("G4").Value = 999
$ws.Range("H4").Value = 81786

$ws.Range("D10").Value = 145617
$ws.Range("E10").Value = 19298

$ws.Range("D23").Value = 26800
$ws.Range("E23").Value = 1699

$ws.Range("B52").Value = 8132
$ws.Range("C52").Value = 27
$ws.Range("E52").Value = 7876

$ws.Range("B56").Value = 6278
$ws.Range("C56").Value = 244
$ws.Range("E56").Value = 4127

$ws.Range("B85").Value = 1730
$ws.Range("C85").Value = 30
$ws.Range("D85").Value = 818
$ws.Range("E85").Value = 891

$ws.Range("A87").Value = "Sudan"
$ws.Range("B87").Value = 1526
$ws.Range("C87").Value = 161
$ws.Range("D87").Value = 162
$ws.Range("E87").Value = 1290
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 4
$ws.Range("H87").Value = 74

$ws.Range("A88").Value = "Nueva Zelanda"
$ws.Range("B88").Value = 1497
$ws.Range("C88").Value = 3
$ws.Range("D88").Value = 1386
$ws.Range("E88").Value = 90
$ws.Range("F88").Value = 2
$ws.Range("H88").Value = 21

$ws.Range("A89").Value = "Lituania"
$ws.Range("B89").Value = 1485
$ws.Range("C89").Value = 6
$ws.Range("D89").Value = 833
$ws.Range("E89").Value = 602
$ws.Range("F89").Value = 17
$ws.Range("H89").Value = 50

$ws.Range("A90").Value = "Eslovenia"
$ws.Range("B90").Value = 1460
$ws.Range("C90").Value = 3
$ws.Range("D90").Value = 256
$ws.Range("E90").Value = 1102
$ws.Range("F90").Value = 10
$ws.Range("H90").Value = 102

$ws.Range("A91").Value = "Eslovaquia"
$ws.Range("B91").Value = 1457
$ws.Range("D91").Value = 959
$ws.Range("E91").Value = 472
$ws.Range("F91").Value = 5
$ws.Range("H91").Value = 26

$ws.Range("B128").Value = 422
$ws.Range("C128").Value = 8
$ws.Range("D128").Value = 205
$ws.Range("E128").Value = 207

$ws.Range("A192").Value = "Belice"
$ws.Range("D192").Value = 16
$ws.Range("H192").Value = 2

$ws.Range("A193").Value = "Nueva Caledonia"
$ws.Range("D193").Value = 18
$ws.Range("H193").Value = 0

$ws.Range("A198").Value = "Dominica"
$ws.Range("D198").Value = 15
$ws.Range("H198").Value = 0

$ws.Range("A199").Value = "Curazao"
$ws.Range("D199").Value = 14
$ws.Range("H199").Value = 1

